# "Tried to implement Penality Reward System (unfinished)"
#
# Weekly Quantity sheet: drop the two weeks that got folded into the
# penalty/reward recompute (old rows 19 & 20), and correct the
# already-adjusted quantity for the 2023-06-11 week (row 14).
#
# Monthly Trend sheet: correct the rolled-up quantities for the two
# affected months (rows 7 & 8).

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B14").Value = 880
$wsWeekly.Rows("19:20").Delete()

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B7").Value = 1180
$wsMonthly.Range("B8").Value = 1040
